$d = $word.ActiveDocument
$d.Bookmarks.ShowHidden = $true

# ------------------------------------------------------------------
# 1) "Relaciones con la H" -> "Relaciones con la h" (lowercase h),
#    bold the whole heading, and move the reserved "_GoBack" bookmark
#    so it sits right between the two runs (i.e. right after "la h").
#    "_GoBack" is a singleton bookmark: re-adding it under this name
#    relocates it, automatically removing it from its old location
#    further down in the document.
# ------------------------------------------------------------------
$rng = $d.Range(0, 0)
if ($rng.Find.Execute("Relaciones con la H")) {
    $rng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng)

    # Lower-case just the last letter ("H" -> "h"). The bookmark we
    # just added already splits the run at this position, so editing
    # this single trailing character keeps the two runs separate
    # instead of Word re-merging them back into one run.
    $hChar = $d.Range($rng.End - 1, $rng.End)
    $hChar.Text = "h"
}

$rngBold1 = $d.Range(0, 0)
if ($rngBold1.Find.Execute("Relaciones con la h")) {
    $d.Paragraphs(2).Range.Bold = 1
}

# ------------------------------------------------------------------
# 2) "Estándar" -> bold
# ------------------------------------------------------------------
$rngEst = $d.Range(0, 0)
if ($rngEst.Find.Execute("Estándar")) {
    $rngEst.Paragraphs(1).Range.Bold = 1
}

# ------------------------------------------------------------------
# 3) "Competencias" -> bold + remove the attached comment
# ------------------------------------------------------------------
$rngComp = $d.Range(0, 0)
if ($rngComp.Find.Execute("Competencias")) {
    $rngComp.Paragraphs(1).Range.Bold = 1
}

for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments($i).Delete()
}

# ------------------------------------------------------------------
# 4) "Estrategia didáctica" -> bold
# ------------------------------------------------------------------
$rngEstr = $d.Range(0, 0)
if ($rngEstr.Find.Execute("Estrategia didáctica")) {
    $rngEstr.Paragraphs(1).Range.Bold = 1
}
